$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3964.2856
$ws.Range("H77").Value = 3964.2856
$ws.Range("H113").Value = 62503996
$ws.Range("I113").Value = 125003740
$ws.Range("K113").Value = 125003740
$ws.Range("M113").Value = -125000486
$ws.Range("H132").Value = 4770.6855
$ws.Range("I132").Value = 4436.75
$ws.Range("K132").Value = 13310.25
$ws.Range("M132").Value = -10780.25
$ws.Range("H137").Value = 2292.2307
$ws.Range("I137").Value = 2306.8572
$ws.Range("K137").Value = 6920.571599999999
$ws.Range("M137").Value = -4370.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6378.385
$ws.Range("I32").Value = 6593.56
$ws.Range("J32").Value = 999
$ws.Range("K32").Value = 6593.56
$ws.Range("L32").Value = 999
$ws.Range("M32").Value = -6306.56
$ws.Range("N32").Value = -1573
$ws.Range("H63").Value = 70593304
$ws.Range("I63").Value = 200002180
$ws.Range("J63").Value = 16672941
$ws.Range("K63").Value = 200002180
$ws.Range("L63").Value = 16672941
$ws.Range("M63").Value = -200001494
$ws.Range("N63").Value = -16674313
$ws.Range("H66").Value = 70593304
$ws.Range("I66").Value = 200002180
$ws.Range("J66").Value = 16672941
$ws.Range("K66").Value = 1000010900
$ws.Range("L66").Value = 83364705
$ws.Range("M66").Value = -1000007468
$ws.Range("N66").Value = -83371569
$ws.Range("H74").Value = 35717572
$ws.Range("I74").Value = 55557070
$ws.Range("K74").Value = 55557070
$ws.Range("M74").Value = -55556196
$ws.Range("H77").Value = 35717572
$ws.Range("I77").Value = 55557070
$ws.Range("K77").Value = 277785350
$ws.Range("M77").Value = -277780982
$ws.Range("H110").Value = 1652.32
$ws.Range("I110").Value = 1427.7273
$ws.Range("K110").Value = 1427.7273
$ws.Range("M110").Value = 617.2727
$ws.Range("H122").Value = 2696.6667
$ws.Range("I122").Value = 1682.8
$ws.Range("K122").Value = 5048.4
$ws.Range("M122").Value = -2598.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 12652.682
$ws.Range("I20").Value = 17909.6
$ws.Range("K20").Value = 17909.6
$ws.Range("M20").Value = -17662.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 4349.5
$ws.Range("J2").Value = 4349.5
$ws.Range("L2").Value = 4349.5
$ws.Range("N2").Value = -4575.5
$ws.Range("H16").Value = 1118.7894
$ws.Range("I16").Value = 966.375
$ws.Range("K16").Value = 966.375
$ws.Range("M16").Value = -679.375
$ws.Range("H31").Value = 3230.7778
$ws.Range("I31").Value = 2328.8235
$ws.Range("J31").Value = 3509.5637
$ws.Range("K31").Value = 2328.8235
$ws.Range("L31").Value = 3509.5637
$ws.Range("M31").Value = -2033.8235
$ws.Range("N31").Value = -4099.563700000001
$ws.Range("H34").Value = 3230.7778
$ws.Range("I34").Value = 2328.8235
$ws.Range("J34").Value = 3509.5637
$ws.Range("K34").Value = 2328.8235
$ws.Range("L34").Value = 3509.5637
$ws.Range("M34").Value = -2126.8235
$ws.Range("N34").Value = -3913.5637
$ws.Range("H113").Value = 1118.7894
$ws.Range("I113").Value = 966.375
$ws.Range("K113").Value = 966.375
$ws.Range("M113").Value = 1203.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4004
$ws.Range("J5").Value = 4286.273
$ws.Range("L5").Value = 12858.819
$ws.Range("N5").Value = -13082.819
$ws.Range("H13").Value = 10187.8
$ws.Range("J13").Value = 106.333336
$ws.Range("L13").Value = 319.000008
$ws.Range("N13").Value = -655.000008
$ws.Range("H135").Value = 4004
$ws.Range("J135").Value = 4286.273
$ws.Range("L135").Value = 38576.457
$ws.Range("N135").Value = -43646.457
$ws.Range("H141").Value = 5945.222
$ws.Range("I141").Value = 6001
$ws.Range("K141").Value = 18003
$ws.Range("M141").Value = -12823

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11873
$ws.Range("J70").Value = 11998.6
$ws.Range("L70").Value = 11998.6
$ws.Range("N70").Value = -12538.6
$ws.Range("H73").Value = 11873
$ws.Range("J73").Value = 11998.6
$ws.Range("L73").Value = 11998.6
$ws.Range("N73").Value = -13870.6
$ws.Range("H113").Value = 2418.577
$ws.Range("I113").Value = 1461.25
$ws.Range("J113").Value = 3239.1428
$ws.Range("K113").Value = 1461.25
$ws.Range("L113").Value = 3239.1428
$ws.Range("M113").Value = 708.75
$ws.Range("N113").Value = -7579.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 10900
$ws.Range("I3").Value = 10900
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10900
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -10788
$ws.Range("N3").ClearContents()
$ws.Range("H15").Value = 10900
$ws.Range("I15").Value = 10900
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 10900
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -10730
$ws.Range("N15").ClearContents()
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H82").Value = 696.4286
$ws.Range("I82").Value = 664.5
$ws.Range("K82").Value = 664.5
$ws.Range("M82").Value = -303.5
$ws.Range("H85").Value = 696.4286
$ws.Range("I85").Value = 664.5
$ws.Range("K85").Value = 664.5
$ws.Range("M85").Value = 583.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 35000
$ws.Range("J31").Value = 35000
$ws.Range("L31").Value = 35000
$ws.Range("N31").Value = -35696
$ws.Range("H81").Value = 11768287
$ws.Range("I81").Value = 2179.7273
$ws.Range("J81").Value = 33339484
$ws.Range("K81").Value = 4359.4546
$ws.Range("L81").Value = 66678968
$ws.Range("M81").Value = -3298.4546
$ws.Range("N81").Value = -66681090
$ws.Range("H84").Value = 11768287
$ws.Range("I84").Value = 2179.7273
$ws.Range("J84").Value = 33339484
$ws.Range("K84").Value = 21797.273
$ws.Range("L84").Value = 333394840
$ws.Range("M84").Value = -16493.273
$ws.Range("N84").Value = -333405448
$ws.Range("H122").Value = 2527.1304
$ws.Range("I122").Value = 2351.8667
$ws.Range("K122").Value = 7055.6001
$ws.Range("M122").Value = -4605.6001
$ws.Range("H126").Value = 1611.5555
$ws.Range("I126").Value = 1321
$ws.Range("J126").Value = 2367
$ws.Range("K126").Value = 3963
$ws.Range("L126").Value = 7101
$ws.Range("M126").Value = -1493
$ws.Range("N126").Value = -12041
